# Add the new "2021" data column (O) to the insurance-indicators table,
# reusing the formatting already applied to the previous year's column (N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: year header
$ws.Range("N3").Copy()
$ws.Range("O3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("O3").Value = 2021

# Row 4: number of reporting insurance companies
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("O4").Value = 14

# Row 5: insurance premiums received
$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("O5").Value = 1252.8

$excel.CutCopyMode = $false

# Match the author's final selection on the worksheet
$ws.Range("O9").Select()
